$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the current row 2 (pushing all existing data rows
# down by 4) to make room for 4 new days of data (17-20 Dec 2021).
$ws.Rows("2:5").Insert()

# The freshly inserted rows default to the formatting of the row above
# (the header row). Copy the date/number formatting used by the rest of
# the data rows (now starting at row 6, the original row 2) down onto the
# new rows.
$ws.Range("A6:B6").Copy()
$ws.Range("A2:B5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows("2:5").RowHeight = $ws.Rows("6").RowHeight

# Populate the 4 new rows with the latest reported NSW second-dose figures.
$ws.Range("A2").Value2 = 44550
$ws.Range("B2").Value2 = 6437227
$ws.Range("A3").Value2 = 44549
$ws.Range("B3").Value2 = 6436676
$ws.Range("A4").Value2 = 44548
$ws.Range("B4").Value2 = 6435466
$ws.Range("A5").Value2 = 44547
$ws.Range("B5").Value2 = 6432466

# Restore the active cell/selection as it was left after the edit.
$ws.Range("C9").Select() | Out-Null
